$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 9 (pushes nothing down since it's the last row; it
# inherits formatting from row 8, same as Excel's native row-insert).
$ws.Rows.Item(9).Insert()

# Row 9 gets the Country/Date that used to belong to row 8, plus the
# OLD C8:H8 values that are being "moved" down a row.
$ws.Range("A9").Value = $ws.Range("A8").Value2
$ws.Range("B9").Value = $ws.Range("B8").Value2

$ws.Range("C9").Value = 0.09630539297713012
$ws.Range("D9").Value = 0.05367771623038212
$ws.Range("E9").Value = -0.001568691807165634
$ws.Range("F9").Value = 0.05524640803754775
$ws.Range("G9").Value = 0.02761047347365702
$ws.Range("H9").Value = 0.9723895265263429

# Row 8 keeps its original Country/Date (A8/B8), but C8:H8 become the
# new values from the diff.
$ws.Range("C8").Value = 0.09217945242613174
$ws.Range("D8").Value = 0.04955177567938374
$ws.Range("E8").Value = -0.004335454964670524
$ws.Range("F8").Value = 0.05388723064405426
$ws.Range("G8").Value = 0.07446332850061535
$ws.Range("H8").Value = 0.9255366714993846
